# Update the "Estado de Cuenta" worker debt table (rows 16-29):
# Eliminates the previous period ordering and re-adds the records with
# periods in descending order (2109 -> 2103), grouping all records for
# JUAN DIEGO ALCALA CABARCAS (CC) first, followed by all records for
# WILLIAM JOSE RAMOS MARIN (PE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker 1: CC 1050970657 - JUAN DIEGO ALCALA CABARCAS
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1050970657"
$ws.Range("D16").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E16").Value = "2109"
$ws.Range("F16").Value = 24227
$ws.Range("G16").Value = 908526

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1050970657"
$ws.Range("D17").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E17").Value = "2108"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1050970657"
$ws.Range("D18").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E18").Value = "2107"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1050970657"
$ws.Range("D19").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E19").Value = "2106"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1050970657"
$ws.Range("D20").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E20").Value = "2105"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 908526

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1050970657"
$ws.Range("D21").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E21").Value = "2104"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1050970657"
$ws.Range("D22").Value = "JUAN DIEGO ALCALA CABARCAS"
$ws.Range("E22").Value = "2103"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

# Worker 2: PE 954073205101971 - WILLIAM JOSE RAMOS MARIN
$ws.Range("B23").Value = "PE"
$ws.Range("C23").Value = "954073205101971"
$ws.Range("D23").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E23").Value = "2109"
$ws.Range("F23").Value = 24227
$ws.Range("G23").Value = 908526

$ws.Range("B24").Value = "PE"
$ws.Range("C24").Value = "954073205101971"
$ws.Range("D24").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E24").Value = "2108"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 908526

$ws.Range("B25").Value = "PE"
$ws.Range("C25").Value = "954073205101971"
$ws.Range("D25").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E25").Value = "2107"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("B26").Value = "PE"
$ws.Range("C26").Value = "954073205101971"
$ws.Range("D26").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E26").Value = "2106"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 908526

$ws.Range("B27").Value = "PE"
$ws.Range("C27").Value = "954073205101971"
$ws.Range("D27").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E27").Value = "2105"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 908526

$ws.Range("B28").Value = "PE"
$ws.Range("C28").Value = "954073205101971"
$ws.Range("D28").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E28").Value = "2104"
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 908526

$ws.Range("B29").Value = "PE"
$ws.Range("C29").Value = "954073205101971"
$ws.Range("D29").Value = "WILLIAM JOSE RAMOS MARIN"
$ws.Range("E29").Value = "2103"
$ws.Range("F29").Value = 26650
$ws.Range("G29").Value = 908526
